# Weekly update: insert a new price record for "Espinaca" (Vega Central
# Mapocho de Santiago) at row 567, pushing the existing rows 567-596 down
# to 568-597.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 567; Excel shifts rows 567..596 down to 568..597
# and carries the row's cell formatting (e.g. the date style on column D).
$ws.Rows.Item(567).Insert()

# Populate the new row 567 with this week's record.
$ws.Cells.Item(567, 1).Value = 9
$ws.Cells.Item(567, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(567, 3).Value = "Metropolitana"
$ws.Cells.Item(567, 4).Value = 45041
$ws.Cells.Item(567, 5).Value = 13
$ws.Cells.Item(567, 6).Value = 100112012
$ws.Cells.Item(567, 7).Value = "Espinaca"
$ws.Cells.Item(567, 8).Value = "Sin especificar"
$ws.Cells.Item(567, 9).Value = "Primera"
$ws.Cells.Item(567, 10).Value = 160
$ws.Cells.Item(567, 11).Value = 10000
$ws.Cells.Item(567, 12).Value = 11000
$ws.Cells.Item(567, 13).Value = 10500
$ws.Cells.Item(567, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(567, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(567, 16).Value = 1050
$ws.Cells.Item(567, 17).Value = 10
$ws.Cells.Item(567, 18).Value = "Hortaliza"
